# Self perception Inventory.xlsx - fill in the Questionnaire answers
# and move the user's focus/selection over to the "Grille d'évaluation" sheet.

$wb = $excel.ActiveWorkbook
$wsQ = $wb.Worksheets.Item("Questionnaire")
$wsG = $wb.Worksheets.Item("Grille d'évaluation")

# ---- Enter the questionnaire scores (column B) ----
$wsQ.Range("B4").Value  = 0
$wsQ.Range("B5").Value  = 0
$wsQ.Range("B6").Value  = 2
$wsQ.Range("B8").Value  = 2
$wsQ.Range("B9").Value  = 2
$wsQ.Range("B10").Value = 2
$wsQ.Range("B11").Value = 2

$wsQ.Range("B15").Value = 1
$wsQ.Range("B16").Value = 1
$wsQ.Range("B17").Value = 1
$wsQ.Range("B18").Value = 4
$wsQ.Range("B19").Value = 1
$wsQ.Range("B20").Value = 0
$wsQ.Range("B21").Value = 2
$wsQ.Range("B22").Value = 0

$wsQ.Range("B26").Value = 2
$wsQ.Range("B27").Value = 1
$wsQ.Range("B28").Value = 2
$wsQ.Range("B29").Value = 1
$wsQ.Range("B30").Value = 1
$wsQ.Range("B31").Value = 0
$wsQ.Range("B32").Value = 2
$wsQ.Range("B33").Value = 1

$wsQ.Range("B37").Value = 0
$wsQ.Range("B38").Value = 2
$wsQ.Range("B39").Value = 3
$wsQ.Range("B40").Value = 0
$wsQ.Range("B41").Value = 0
$wsQ.Range("B42").Value = 3
$wsQ.Range("B43").Value = 0
$wsQ.Range("B44").Value = 2

$wsQ.Range("B48").Value = 3
$wsQ.Range("B49").Value = 3
$wsQ.Range("B50").Value = 0
$wsQ.Range("B51").Value = 2
$wsQ.Range("B52").Value = 0
$wsQ.Range("B53").Value = 0
$wsQ.Range("B54").Value = 0
$wsQ.Range("B55").Value = 2

$wsQ.Range("B59").Value = 2
$wsQ.Range("B60").Value = 2
$wsQ.Range("B61").Value = 1
$wsQ.Range("B62").Value = 0
$wsQ.Range("B63").Value = 2
$wsQ.Range("B64").Value = 1
$wsQ.Range("B65").Value = 2
$wsQ.Range("B66").Value = 0

$wsQ.Range("B70").Value = 3
$wsQ.Range("B71").Value = 1
$wsQ.Range("B72").Value = 1
$wsQ.Range("B73").Value = 0
$wsQ.Range("B74").Value = 3
$wsQ.Range("B75").Value = 2
$wsQ.Range("B76").Value = 0
$wsQ.Range("B77").Value = 0

$excel.CalculateFull()

# ---- Update sheet selections / which sheet is active ----
$wsQ.Activate()
$wsQ.Range("B76").Select()

$wsG.Activate()
$wsG.Range("P4:Q4").Select()

Write-Host "Questionnaire scores entered; Grille d'évaluation is now active."
